$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.526267290115356
$ws.Range("C1").Value = 1.857581853866577
$ws.Range("D1").Value = 1.444078683853149
$ws.Range("E1").Value = 1.305597305297852
